# Change the "Domain" value on slide 1 from "Twitter" to "Social Networks".
# The shape is "Rectangle 6" (shape index 4 on slide 1), whose text is
# "Domain\t: Twitter" held in two separate runs: "Domain\t: " and "Twitter".
# We only want to replace the second run's text, leaving its formatting
# (lang="en-GB" b="1") and the first run untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange

$fullText = $tr.Text
$target = "Twitter"
$idx = $fullText.IndexOf($target)

if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $target.Length)
    $sub.Text = "Social Networks"
}
